$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textRefs = @("D4", "D5", "D7", "D8", "D9", "D10", "D11", "D12", "D14", "D15", "D16", "D17", "D18", "D20", "D23", "D25", "D27", "D28", "D30", "D31", "D32", "D33", "D34", "D35", "D37", "D39", "D40", "D41", "D42", "D43", "D45", "D46", "D47", "D48", "D49")
foreach ($ref in $textRefs) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = "27.461.55"
$ws.Range("E2").Value = "  -1.29%  "
$ws.Range("D3").Value = "1.832.56"
$ws.Range("E3").Value = "  -1.96%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.75%  "
$ws.Range("D5").Value = "331.19"
$ws.Range("E6").Value = "  -0.72%  "
$ws.Range("D7").Value = "0.4604"
$ws.Range("E7").Value = "  -2.85%  "
$ws.Range("D8").Value = "0.3835"
$ws.Range("E8").Value = "  -1.99%  "
$ws.Range("D9").Value = "46.65"
$ws.Range("E9").Value = "  -0.24%  "
$ws.Range("D10").Value = "0.07884"
$ws.Range("E10").Value = "  -1.14%  "
$ws.Range("D11").Value = "0.9721"
$ws.Range("E11").Value = "  -3.41%  "
$ws.Range("D12").Value = "21.11"
$ws.Range("E12").Value = "  -2.32%  "
$ws.Range("D13").Value = "1.837.27"
$ws.Range("E13").Value = "  -2.27%  "
$ws.Range("D14").Value = "5.887"
$ws.Range("E14").Value = "  -1.68%  "
$ws.Range("D15").Value = "7.060"
$ws.Range("E15").Value = "  -1.30%  "
$ws.Range("D16").Value = "1.001"
$ws.Range("E16").Value = "  -1.02%  "
$ws.Range("D17").Value = "88.05"
$ws.Range("E17").Value = "  -0.05%  "
$ws.Range("D18").Value = "0.06627"
$ws.Range("E18").Value = "  -1.05%  "
$ws.Range("E19").Value = "  -0.94%  "
$ws.Range("D20").Value = "17.12"
$ws.Range("E20").Value = "  +0.96%  "
$ws.Range("E21").Value = "  -0.65%  "
$ws.Range("D22").Value = "27.471.92"
$ws.Range("E22").Value = "  -1.33%  "
$ws.Range("D23").Value = "5.345"
$ws.Range("E23").Value = "  -2.39%  "
$ws.Range("E24").Value = "  -0.87%  "
$ws.Range("D25").Value = "2.305"
$ws.Range("E25").Value = "  -0.97%  "
$ws.Range("D26").Value = "2.064.07"
$ws.Range("E26").Value = "  -1.70%  "
$ws.Range("D27").Value = "157.25"
$ws.Range("E27").Value = "  -0.43%  "
$ws.Range("D28").Value = "19.40"
$ws.Range("E28").Value = "  -1.45%  "
$ws.Range("E29").Value = "  -1.00%  "
$ws.Range("D30").Value = "5.286"
$ws.Range("E30").Value = "  -2.36%  "
$ws.Range("D31").Value = "118.83"
$ws.Range("E31").Value = "  -1.90%  "
$ws.Range("D32").Value = "0.9558"
$ws.Range("E32").Value = "  -1.16%  "
$ws.Range("D33").Value = "0.09287"
$ws.Range("E33").Value = "  -1.92%  "
$ws.Range("D34").Value = "3.577"
$ws.Range("D35").Value = "5.245"
$ws.Range("E35").Value = "  -1.04%  "
$ws.Range("E36").Value = "  -2.06%  "
$ws.Range("D37").Value = "0.05942"
$ws.Range("E37").Value = "  -1.53%  "
$ws.Range("E38").Value = "  -0.98%  "
$ws.Range("D39").Value = "8.061"
$ws.Range("E39").Value = "  -0.86%  "
$ws.Range("D40").Value = "1.153"
$ws.Range("E40").Value = "  -4.01%  "
$ws.Range("D41").Value = "0.5804"
$ws.Range("E41").Value = "  -2.03%  "
$ws.Range("D42").Value = "0.1842"
$ws.Range("E42").Value = "  -2.36%  "
$ws.Range("D43").Value = "10.02"
$ws.Range("E43").Value = "  -2.38%  "
$ws.Range("E44").Value = "  +2.42%  "
$ws.Range("D45").Value = "0.5491"
$ws.Range("E45").Value = "  -2.47%  "
$ws.Range("D46").Value = "12.03"
$ws.Range("E46").Value = "  -0.27%  "
$ws.Range("D47").Value = "1.872"
$ws.Range("E47").Value = "  -2.22%  "
$ws.Range("D48").Value = "0.06652"
$ws.Range("E48").Value = "  -1.81%  "
$ws.Range("D49").Value = "110.45"
$ws.Range("E49").Value = "  -1.24%  "
$ws.Range("E50").Value = "  -2.17%  "
$ws.Range("E51").Value = "  -0.82%  "

foreach ($ref in $textRefs) {
    $ws.Range($ref).Style = "Normal"
}

Write-Output "Applied cryptos update"
